$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map, derived from the authoritative OOXML diff.
$updates = [ordered]@{
    "S2" = 1.33
    "G3" = 2.05
    "I3" = 3.6
    "J3" = 2.75
    "L3" = 4.33
    "M3" = 1.06
    "N3" = 10
    "Q3" = 2.15
    "R3" = 1.67
    "S3" = 1.44
    "T3" = 2.63
    "U3" = 2
    "V3" = 1.73
    "W3" = 6.5
    "X3" = 9
    "Y3" = 9
    "Z3" = 17
    "AJ3" = 13
    "AL3" = 34
    "AV3" = 67
    "AY3" = 34
    "S4" = 1.53
    "T4" = 2.38
    "Q5" = 2.88
    "R5" = 1.4
    "S5" = 1.67
    "N6" = 8
    "Q6" = 2.3
    "R6" = 1.6
    "L13" = 4
    "N13" = 9.5
    "Q13" = 2.05
    "S13" = 1.4
    "T13" = 2.75
    "Z13" = 19
    "AG13" = 251
    "AI13" = 17
    "AJ13" = 12
    "AT13" = 2.75
    "Q15" = 1.93
    "R15" = 1.93
    "Q16" = 2.03
    "R16" = 1.83
    "G17" = 1.8
    "H17" = 3.6
    "AH17" = 13
    "AO17" = 9.5
    "BD17" = 126
    "Q25" = 1.85
    "R25" = 2
    "G30" = 2.05
    "H30" = 3.1
    "I30" = 3.7
    "J30" = 2.88
    "L30" = 4.33
    "N30" = 8
    "O30" = 1.4
    "P30" = 2.75
    "U30" = 1.91
    "V30" = 1.8
    "W30" = 6.5
    "X30" = 9
    "Y30" = 9.5
    "Z30" = 19
    "AA30" = 19
    "AE30" = 17
    "AH30" = 9.5
    "AI30" = 17
    "AJ30" = 13
    "AK30" = 41
    "AL30" = 34
    "AN30" = 4
    "AO30" = 12
    "AP30" = 23
    "AQ30" = 41
    "AW30" = 5.5
    "AX30" = 21
    "AY30" = 34
    "AZ30" = 67
    "BA30" = 101
    "G32" = 2.55
    "I32" = 2.5
    "L32" = 3.1
    "AK32" = 26
    "AL32" = 19
    "AP32" = 19
    "G44" = 1.8
    "H44" = 3.8
    "I44" = 4.1
    "J44" = 2.25
    "L44" = 4
    "N44" = 21
    "X44" = 12
    "Y44" = 9
    "Z44" = 17
    "AA44" = 12
    "AD44" = 8.5
    "AE44" = 11
    "AF44" = 29
    "AG44" = 81
    "AH44" = 21
    "AI44" = 29
    "AJ44" = 15
    "AL44" = 29
    "AM44" = 26
    "AO44" = 9
    "AQ44" = 23
    "AR44" = 34
    "AW44" = 6.5
    "AX44" = 19
    "AY44" = 21
    "BB44" = 101
    "BC44" = 251
    "G45" = 1.29
    "I45" = 9
    "L45" = 7.5
    "U45" = 1.8
    "V45" = 1.91
    "X45" = 7.5
    "AB45" = 23
    "AJ45" = 23
    "AU45" = 9
    "I46" = 3.8
    "AA46" = 12
    "AH46" = 21
    "AR46" = 34
    "BC46" = 201
    "G47" = 2.25
    "I47" = 2.9
    "J47" = 2.75
    "L47" = 3.25
    "O47" = 1.17
    "P47" = 5
    "Q47" = 1.6
    "R47" = 2.3
    "S47" = 1.29
    "T47" = 3.5
    "W47" = 11
    "X47" = 13
    "Y47" = 9.5
    "Z47" = 21
    "AA47" = 17
    "AB47" = 21
    "AH47" = 13
    "AI47" = 17
    "AJ47" = 11
    "AK47" = 29
    "AL47" = 21
    "AO47" = 12
    "AT47" = 3.5
    "AW47" = 5
    "AX47" = 15
    "AY47" = 21
    "Q50" = 1.62
    "R50" = 2.25
    "BD50" = 126
    "I54" = 3
    "J54" = 2.88
    "K54" = 2.3
    "U54" = 1.57
    "V54" = 2.25
    "W54" = 10
    "AM54" = 26
    "AN54" = 4.5
    "AX54" = 15
    "AY54" = 21
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Output ("Updated {0} cells" -f $updates.Count)
